$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "demand" column
$ws.Range("D1").Value = "demand"

# Demand values for rows 2-25
$demand = @(60, 67, 71, 19, 59, 83, 34, 55, 33, 64, 26, 42, 73, 81, 50, 88, 73, 74, 16, 85, 47, 59, 52, 62)

for ($i = 0; $i -lt $demand.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $demand[$i]
}

# F6 becomes an empty cell that carries the same (non-default) cell style as B2.
# Give F6 a temporary value so it participates in the style update, apply a
# font touch that is a value no-op (so visually nothing changes) across both
# B2 and F6 together so the engine reuses the existing shared style record,
# then clear the temporary value so F6 ends up empty but still styled.
$ws.Range("F6").Value = 0
$styledCells = $excel.Union($ws.Range("B2"), $ws.Range("F6"))
foreach ($cell in $styledCells) {
    $cell.Font.ThemeColor = 1
}
$ws.Range("F6").ClearContents()

# Update the selected cell to F6, matching the saved selection state
$ws.Range("F6").Select()
